# Remove the post "「車のバッテリーが要る！」" (row 199) from the sheet.
# Deleting the entire row shifts all subsequent rows up by one, matching
# the target diff (row 199 gone, rows 200-263 renumbered to 199-262, and
# the sheet dimension shrinking from A1:C263 to A1:C262).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(199).Delete()
